# Resize/reposition the two "{{画像_N}}" placeholder textboxes on slide 1
# and switch their text frames from top-left/spAutoFit to a centered,
# non-autofitting layout (wrap stays "none").
#
# Point <-> EMU note: Shape.Left/Top/Width/Height are expressed in points
# (1 pt = 12700 EMU) and are stored internally as 32-bit floats before being
# converted back to EMU, so the literals below are chosen (one ULP above the
# naive `emu/12700` value where needed) so the round-trip lands exactly on
# the target EMU values from the target OOXML.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# TextBox 4 (id=5, "{{画像_1}}") -> off 284921,1614186 ext 2132970,3079561
$sh1 = $s.Shapes.Item(3)
$sh1.Left = 22.434724807739258
$sh1.Top = 127.10126495361328
$sh1.Width = 167.95040893554688
$sh1.Height = 242.48512268066406
$sh1.TextFrame.WordWrap = [Microsoft.Office.Core.MsoTriState]::msoFalse
$sh1.TextFrame.AutoSize = [Microsoft.Office.Core.MsoAutoSize]::msoAutoSizeNone
$sh1.TextFrame.VerticalAnchor = [Microsoft.Office.Core.MsoVerticalAnchor]::msoAnchorMiddle

# TextBox 7 (id=8, "{{画像_2}}") -> off 284921,5853858 ext 2132970,3079561
$sh2 = $s.Shapes.Item(6)
$sh2.Left = 22.434724807739258
$sh2.Top = 460.9337158203125
$sh2.Width = 167.95040893554688
$sh2.Height = 242.48512268066406
$sh2.TextFrame.WordWrap = [Microsoft.Office.Core.MsoTriState]::msoFalse
$sh2.TextFrame.AutoSize = [Microsoft.Office.Core.MsoAutoSize]::msoAutoSizeNone
$sh2.TextFrame.VerticalAnchor = [Microsoft.Office.Core.MsoVerticalAnchor]::msoAnchorMiddle
